# Updates the cryptos price/volume table to the latest scraped values.
# Price values that look numeric (e.g. "213.50") must stay stored as TEXT
# (matching the sheet's existing inline-string cells), so those are written
# via a NumberFormat="@" / Style="Normal" round-trip that forces text
# interpretation without leaving a residual style on the cell. Values that
# are already unambiguous text (multiple dots, subscript digits, etc.) or
# the percentage column (which always carries padding spaces) are set
# directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.434.26"
$ws.Range("E2").Value = "  -0.28%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.584.23"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "213.50"
$ws.Range("E5").Value = "  +0.27%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.29%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - OKB
Set-TextValue "D8" "44.66"
$ws.Range("E8").Value = "  -0.34%  "

# Row 9 - Solana
Set-TextValue "D9" "23.93"
$ws.Range("E9").Value = "  -1.17%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.248"
$ws.Range("E10").Value = "  -1.62%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0590"
$ws.Range("E11").Value = "  -1.73%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.98%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "1.810.70"
$ws.Range("E13").Value = "  -0.20%  "

# Row 14 - WrappedEther
Set-TextValue "D14" "1.588.15"
$ws.Range("E14").Value = "  +0.03%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -0.75%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  -1.80%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "28.443.17"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18 - Litecoin
Set-TextValue "D18" "62.12"
$ws.Range("E18").Value = "  -1.48%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "230.15"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.50%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0689"
$ws.Range("E21").Value = "  -2.38%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.03%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -3.26%  "

# Row 24 - Avalanche
Set-TextValue "D24" "9.16"
$ws.Range("E24").Value = "  -1.83%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.84%  "

# Row 26 - Monero
Set-TextValue "D26" "151.76"
$ws.Range("E26").Value = "  -0.02%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "15.05"
$ws.Range("E27").Value = "  -1.11%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -1.64%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -1.85%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  -0.01%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +2.68%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -1.40%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -1.27%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -2.41%  "

# Row 35 - Maker
Set-TextValue "D35" "1.395.49"
$ws.Range("E35").Value = "  +0.22%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +6.86%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -4.85%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  +0.21%  "

# Row 39 - MXToken
Set-TextValue "D39" "2.65"
$ws.Range("E39").Value = "  +0.93%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.76%  "

# Row 41 - ImmutableX
Set-TextValue "D41" "0.523"
$ws.Range("E41").Value = "  -3.24%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.03%  "

# Row 43 - ARBITRUM
Set-TextValue "D43" "0.792"
$ws.Range("E43").Value = "  -2.39%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  +0.90%  "

# Row 45 - FraxShare
Set-TextValue "D45" "5.45"
$ws.Range("E45").Value = "  -3.58%  "

# Row 46 - Kaspa
Set-TextValue "D46" "0.0458"
$ws.Range("E46").Value = "  -1.36%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  -2.27%  "

# Row 48 - Aave
Set-TextValue "D48" "62.88"
$ws.Range("E48").Value = "  -0.03%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "1.722.00"
$ws.Range("E49").Value = "  -0.07%  "

# Row 50 - Quant
Set-TextValue "D50" "86.68"
$ws.Range("E50").Value = "  -0.29%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  -1.35%  "
